$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.922.57'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.880.78'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.80'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4604'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3875'
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07862'
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9869'
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.78'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.828.28'
$ws.Range("E12").Value = '  -2.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.992'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.649'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06960'
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.19'
$ws.Range("E16").Value = '  +0.67%  '
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009964'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.98'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("E20").Value = '  -0.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.915.88'
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("E22").Value = '  -1.85%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.087'
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.16'
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.009'
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.64'
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.922'
$ws.Range("E29").Value = '  -2.29%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9021'
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.266'
$ws.Range("E32").Value = '  -1.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.316'
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.256'
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.180'
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05743'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02072'
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.626'
$ws.Range("E39").Value = '  -4.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5646'
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1763'
$ws.Range("E41").Value = '  -1.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.678'
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.273'
$ws.Range("E43").Value = '  +5.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.87'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5347'
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07042'
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.843'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '112.72'
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.532'
$ws.Range("E49").Value = '  +1.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.068'
$ws.Range("E50").Value = '  -4.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.65'
$ws.Range("E51").Value = '  -0.26%  '
